$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "J3RWJ6MRC7YZ"
$ws.Range("A3").Value = "JYKGY5W23KHC"
$ws.Range("A4").Value = "X00J9QE1TGY2"
$ws.Range("A11").Value = "7VADX0RWN3NT"
$ws.Range("A12").Value = "AA4GDR1CGGRM"
$ws.Range("A13").Value = "VB8SVSWCP89Z"
$ws.Range("A14").Value = "20P7SKTWSCG7"
$ws.Range("A15").Value = "2578DVHZ2NY8"
$ws.Range("A16").Value = "05TDN56XHKZ0"
$ws.Range("A17").Value = "3WKDV7ZXQN1Y"

$ws.Range("A2:A4").Select()
